$p = $ppt.ActivePresentation

# Slide 1 shape order: Title 1, Picture 1, TextBox 3, Content Placeholder 3
# Slide 1: Title "Slide" + " " + "1" -> single run "Slide 1"
$s1 = $p.Slides.Item(1)
$titleTr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$titleTr1.Text = " "
$titleTr1.Text = "Slide 1"
# Slide 1: Caption "an" + " " + "image" -> single run "an image"
$capTr1 = $s1.Shapes.Item(3).TextFrame.TextRange
$capTr1.Text = " "
$capTr1.Text = "an image"

# Slide 2 shape order: Title 1, Content Placeholder 2, Picture 1, TextBox 3
# Slide 2: Title "Slide" + " " + "2" -> single run "Slide 2"
$s2 = $p.Slides.Item(2)
$titleTr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$titleTr2.Text = " "
$titleTr2.Text = "Slide 2"
# Slide 2: Caption "an" + " " + "image" -> single run "an image"
$capTr2 = $s2.Shapes.Item(4).TextFrame.TextRange
$capTr2.Text = " "
$capTr2.Text = "an image"
